# Dietician_testdata.xlsx edit script
# Adds "additional details" scenario columns (Middlename, SecondaryContact)
# plus two new "create_dietician_AdditionalData" test rows on the
# Dietician_Create sheet, and fixes a pincode typo.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dietician_Create")

# ---------------------------------------------------------------------
# 1) Capture the existing K/L ("scenario" / error-message) column values
#    for rows 1-9 before we touch anything, then remove those two
#    columns outright so we can rebuild K/L (new Middlename /
#    SecondaryContact headers) and M/N (the relocated scenario columns)
#    cleanly without leaving stray empty cells behind.
# ---------------------------------------------------------------------
$colKVals = @{}
$colLVals = @{}
for ($r = 1; $r -le 9; $r++) {
    $colKVals[$r] = $ws.Cells.Item($r, 11).Value2
    $colLVals[$r] = $ws.Cells.Item($r, 12).Value2
}

$ws.Columns.Item(11).EntireColumn.Delete() | Out-Null
$ws.Columns.Item(11).EntireColumn.Delete() | Out-Null

# ---------------------------------------------------------------------
# 2) Fix the HospitalPincode typo on row 6 (60050 -> 600050)
# ---------------------------------------------------------------------
$ws.Range("H6").Value = 600050

# ---------------------------------------------------------------------
# 3) Re-create the old "scenario" (K) / message (L) values in the new
#    M / N columns for rows 1-9, restoring wrap-text formatting.
# ---------------------------------------------------------------------
for ($r = 1; $r -le 9; $r++) {
    $ws.Cells.Item($r, 13).Value = $colKVals[$r]
    $ws.Cells.Item($r, 14).Value = $colLVals[$r]
    $ws.Cells.Item($r, 14).WrapText = $true
}
# Row 1 ("scenario" header) keeps the original unformatted look; the
# rest (M2:M9) get wrap-text like the rest of the header/data cells.
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 13).WrapText = $true
}
$ws.Cells.Item(1, 14).WrapText = $true

# ---------------------------------------------------------------------
# 4) New Middlename / SecondaryContact header cells (K1 / L1)
# ---------------------------------------------------------------------
$ws.Range("K1").Value = "Middlename"
$ws.Range("L1").Value = "SecondaryContact"
$ws.Range("K1").WrapText = $true
$ws.Range("L1").WrapText = $true

# ---------------------------------------------------------------------
# 5) Row heights: data rows 2-9 shrink from 48 to 32 (header row 1
#    already is 32 and stays that way).
# ---------------------------------------------------------------------
for ($r = 2; $r -le 9; $r++) {
    $ws.Rows.Item($r).RowHeight = 32
}

# ---------------------------------------------------------------------
# 6) New scenario rows 10 & 11 - "create_dietician_AdditionalData"
# ---------------------------------------------------------------------
$ws.Range("K10").Value = "numpy"
$ws.Range("L10").Value = 8456798123
$ws.Range("M10").Value = "create_dietician_AdditionalData"
$ws.Range("N10").Value = "valid_additional_data"

$ws.Range("K11").Value = 123
$ws.Range("L11").Value = "7123661a58"
$ws.Range("M11").Value = "create_dietician_AdditionalData"
$ws.Range("N11").Value = "invalid_additional_data"

$ws.Range("K10:N11").WrapText = $true
$ws.Rows.Item(10).RowHeight = 32
$ws.Rows.Item(11).RowHeight = 32

# ---------------------------------------------------------------------
# 7) Column widths for the new L / M / N columns
# ---------------------------------------------------------------------
$ws.Columns.Item(12).ColumnWidth = 10.25
$ws.Columns.Item(13).ColumnWidth = 18.916666666666668
$ws.Columns.Item(14).ColumnWidth = 22.416666666666668

# ---------------------------------------------------------------------
# 8) Selection / active cell moves to H6 (where the pincode fix was
#    made) instead of the old C13.
# ---------------------------------------------------------------------
$ws.Range("H6").Select() | Out-Null
